$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price values in column D (Price).
# The source data stores prices as text, so force a Text number format
# before assigning the value to avoid Excel auto-converting it to a number
# (which would lose significant trailing/leading zeros, e.g. 0.001620 -> 0.00162).

$priceCells = @("D2", "D4", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D28", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "246.68"
$ws.Range("D4").Value = "5.293"
$ws.Range("D6").Value = "3.394"
$ws.Range("D7").Value = "6.375"
$ws.Range("D8").Value = "0.8158"
$ws.Range("D9").Value = "0.9659"
$ws.Range("D11").Value = "0.07420"
$ws.Range("D12").Value = "0.03150"
$ws.Range("D13").Value = "0.03047"
$ws.Range("D14").Value = "0.09284"
$ws.Range("D15").Value = "3.587"
$ws.Range("D16").Value = "0.001620"
$ws.Range("D17").Value = "0.04696"
$ws.Range("D18").Value = "0.0005765"
$ws.Range("D19").Value = "0.006445"
$ws.Range("D20").Value = "0.005053"
$ws.Range("D22").Value = "0.0001498"
$ws.Range("D23").Value = "3.759"
$ws.Range("D24").Value = "2.125"
$ws.Range("D25").Value = "0.3253"
$ws.Range("D26").Value = "0.1248"
$ws.Range("D28").Value = "0.0003097"
$ws.Range("D40").Value = "0.03930"
$ws.Range("D41").Value = "0.007028"
$ws.Range("D42").Value = "0.1045"
$ws.Range("D43").Value = "0.003397"
$ws.Range("D44").Value = "0.007807"
$ws.Range("D45").Value = "0.00005806"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("D47").Value = "0.0005495"
$ws.Range("D48").Value = "0.6793"
$ws.Range("D49").Value = "0.1522"
$ws.Range("D50").Value = "0.00002098"
